$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.954.13"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3
$ws.Range("D3").Value = "1.554.65"
$ws.Range("E3").Value = "  +0.22%  "

# Row 4
$ws.Range("E4").Value = "  -0.58%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.07"
$ws.Range("E5").Value = "  +0.03%  "

# Row 6
$ws.Range("E6").Value = "  +0.67%  "

# Row 7
$ws.Range("E7").Value = "  -0.53%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.06"
$ws.Range("E8").Value = "  +2.18%  "

# Row 9
$ws.Range("E9").Value = "  -0.22%  "

# Row 10
$ws.Range("E10").Value = "  +1.38%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0856"
$ws.Range("E11").Value = "  -0.46%  "

# Row 12
$ws.Range("D12").Value = "1.775.54"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13
$ws.Range("D13").Value = "1.557.66"
$ws.Range("E13").Value = "  +0.32%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.76"
$ws.Range("E14").Value = "  +1.23%  "

# Row 15
$ws.Range("E15").Value = "  +1.04%  "

# Row 16
$ws.Range("D16").Value = "26.944.44"
$ws.Range("E16").Value = "  +0.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.72"
$ws.Range("E17").Value = "  -0.15%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0708"
$ws.Range("E18").Value = "  +2.92%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.58"
$ws.Range("E19").Value = "  +1.12%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.32"
$ws.Range("E20").Value = "  +1.23%  "

# Row 21
$ws.Range("E21").Value = "  -0.60%  "

# Row 22
$ws.Range("E22").Value = "  +1.92%  "

# Row 23
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("E24").Value = "  -2.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.44"
$ws.Range("E25").Value = "  +0.44%  "

# Row 26
$ws.Range("E26").Value = "  -0.51%  "

# Row 27
$ws.Range("E27").Value = "  +0.97%  "

# Row 28
$ws.Range("E28").Value = "  +0.96%  "

# Row 29
$ws.Range("E29").Value = "  -0.55%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0469"
$ws.Range("E30").Value = "  +1.27%  "

# Row 31
$ws.Range("E31").Value = "  -0.37%  "

# Row 32
$ws.Range("E32").Value = "  +0.38%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.12"
$ws.Range("E33").Value = "  +4.21%  "

# Row 34
$ws.Range("D34").Value = "1.418.08"
$ws.Range("E34").Value = "  +0.93%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  +3.28%  "

# Row 36
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  +11.63%  "

# Row 37
$ws.Range("E37").Value = "  +0.41%  "

# Row 38
$ws.Range("E38").Value = "  +0.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.530"
$ws.Range("E39").Value = "  +2.03%  "

# Row 40
$ws.Range("E40").Value = "  -0.14%  "

# Row 41
$ws.Range("E41").Value = "  -0.52%  "

# Row 42
$ws.Range("E42").Value = "  +2.60%  "

# Row 43
$ws.Range("E43").Value = "  +2.21%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.97%  "

# Row 45
$ws.Range("E45").Value = "  +1.36%  "

# Row 46
$ws.Range("E46").Value = "  +0.45%  "

# Row 47
$ws.Range("D47").Value = "1.689.44"
$ws.Range("E47").Value = "  +0.18%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.29"
$ws.Range("E48").Value = "  +1.30%  "

# Row 49
$ws.Range("E49").Value = "  +1.34%  "

# Row 50
$ws.Range("E50").Value = "  +3.01%  "

# Row 51
$ws.Range("E51").Value = "  +0.82%  "
